$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename headers under "table  Questions" (row 5): qname/qimportance/qtype -> question/importance/type
$ws.Range("C5").Value = "question"
$ws.Range("D5").Value = "importance"
$ws.Range("E5").Value = "type(单选or多选)"

# Rename headers under "table Answers" (row 9): qoption -> answer, relavent(...) -> binding(...)
$ws.Range("C9").Value = "answer"
$ws.Range("E9").Value = "binding(该选项关联的下个问题)"

# Update the view state to match: scrolled so row 4 is at top, selection on E15
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E15").Select()
